$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.334.66'
$ws.Range("E2").Value = '  +2.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.057.85'
$ws.Range("E3").Value = '  +4.38%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.55'
$ws.Range("E5").Value = '  +1.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.615'
$ws.Range("E6").Value = '  +3.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '57.89'
$ws.Range("E7").Value = '  +7.04%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +3.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.81'
$ws.Range("E10").Value = '  -0.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0761'
$ws.Range("E11").Value = '  +1.91%  '
$ws.Range("E12").Value = '  +3.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.362.72'
$ws.Range("E13").Value = '  +4.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.46'
$ws.Range("E14").Value = '  +4.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.96'
$ws.Range("E15").Value = '  +5.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.776'
$ws.Range("E16").Value = '  +4.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.19'
$ws.Range("E17").Value = '  +3.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.064.16'
$ws.Range("E18").Value = '  +4.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '37.529.58'
$ws.Range("E19").Value = '  +3.40%  '
$ws.Range("E20").Value = '  +17.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.11'
$ws.Range("E21").Value = '  +2.43%  '
$ws.Range("E22").Value = '  +1.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '226.48'
$ws.Range("E23").Value = '  +2.54%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("E25").Value = '  +3.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.39'
$ws.Range("E26").Value = '  +1.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.78'
$ws.Range("E27").Value = '  +2.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.50'
$ws.Range("E28").Value = '  +13.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.85'
$ws.Range("E29").Value = '  +4.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.13'
$ws.Range("E30").Value = '  +2.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.125'
$ws.Range("E31").Value = '  +1.85%  '
$ws.Range("E32").Value = '  +2.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.51'
$ws.Range("E33").Value = '  +4.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0622'
$ws.Range("E34").Value = '  +3.52%  '
$ws.Range("E35").Value = '  +12.37%  '
$ws.Range("E36").Value = '  +6.47%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.43'
$ws.Range("E37").Value = '  +6.80%  '
$ws.Range("B38").Value = 'BinanceUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.12%  '
$ws.Range("E39").Value = '  +0.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.88'
$ws.Range("E40").Value = '  +12.15%  '
$ws.Range("B41").Value = 'FTXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.59'
$ws.Range("E41").Value = '  +29.19%  '
$ws.Range("B42").Value = 'Cronos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0991'
$ws.Range("E42").Value = '  +12.19%  '
$ws.Range("B43").Value = 'HuobiToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.97'
$ws.Range("E43").Value = '  -1.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.17'
$ws.Range("E44").Value = '  +10.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.472.48'
$ws.Range("E45").Value = '  +1.77%  '
$ws.Range("E46").Value = '  +8.20%  '
$ws.Range("E47").Value = '  +4.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.88'
$ws.Range("E48").Value = '  +8.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.02'
$ws.Range("E49").Value = '  +4.07%  '
$ws.Range("E50").Value = '  +7.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.93'
$ws.Range("E51").Value = '  +2.22%  '
